$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "A ‘Messiah’ for the Multitudes, Freed From History’s Bonds"
$ws.Range("B2").Value = "https://www.nytimes.com/2020/12/21/arts/music/handel-messiah-canada-indigenous.html"
$ws.Range("A3").Value = "Sonos launches microphone-free Arc SL soundbar exclusively at Costco"
$ws.Range("B3").Value = "https://www.theverge.com/2020/12/2/22011112/sonos-arc-sl-no-microphone-now-available-costco"
$ws.Range("A4").Value = "Canada bans mass exports of prescription drugs"
$ws.Range("B4").Value = "https://www.bbc.co.uk/news/world-us-canada-55119428"
$ws.Range("A5").Value = "Canada Sixties Scoop: Indigenous survivors map out their stories"
$ws.Range("B5").Value = "https://www.bbc.co.uk/news/av/world-us-canada-55269251"
$ws.Range("A6").Value = "Gatik’s self-driving box trucks to shuttle groceries for Loblaw in Canada"
$ws.Range("B6").Value = "http://techcrunch.com/2020/11/23/gatiks-self-driving-box-trucks-to-shuttle-groceries-for-loblaw-in-canada/"
$ws.Range("A7").Value = "BRIEF-Canada PM Trudeau says first batch Of Pfizer/BioNTech vaccines arrived in Canada - Reuters India"
$ws.Range("B7").Value = "https://in.reuters.com/article/brief-canada-pm-trudeau-says-first-batch-idINL1N2IU01V"
$ws.Range("A8").Value = "BRIEF-Canada PM Trudeau says first batch Of Pfizer/BioNTech vaccines arrived in Canada - Reuters India"
$ws.Range("B8").Value = "https://in.reuters.com/article/brief-canada-pm-trudeau-says-first-batch-idUKL1N2IU01V"
$ws.Range("A9").Value = "Google’s Nest Hub Max smart screen can now make Zoom calls"
$ws.Range("B9").Value = "http://techcrunch.com/2020/12/14/googles-nest-hub-max-smart-screen-can-now-make-zoom-calls/"
$ws.Range("A10").Value = "Lightspeed acquires restaurant software company Upserve for `$430M"
$ws.Range("B10").Value = "http://techcrunch.com/2020/12/01/lightspeed-acquires-upserve/"
$ws.Range("A11").Value = "ServiceNow is acquiring Element AI, the Canadian startup building AI services for enterprises"
$ws.Range("B11").Value = "http://techcrunch.com/2020/11/30/servicenow-is-acquiring-element-ai-the-canadian-startup-building-ai-services-for-enterprises/"
$ws.Range("A12").Value = "Canada Is Latest to Approve Covid-19 Vaccine—and U.S. Is Likely Next"
$ws.Range("B12").Value = "https://gizmodo.com/canada-is-latest-to-approve-covid-19-vaccine-and-u-s-i-1845844914"
$ws.Range("A13").Value = "Canada extends travel restrictions for those entering the country - Reuters Canada"
$ws.Range("B13").Value = "https://ca.reuters.com/article/us-health-coronavirus-canada-travel-idCAKBN2890XG"
$ws.Range("A14").Value = "Amid surging second coronavirus wave, Canada to unveil more spending - Reuters Canada"
$ws.Range("B14").Value = "https://ca.reuters.com/article/canada-budget-idCAKBN28A1AH"
$ws.Range("A15").Value = "Defense grilling of Canada police witness in Huawei CFOs U.S. extradition case continues - Reuters Canada"
$ws.Range("B15").Value = "https://ca.reuters.com/article/us-usa-huawei-tech-canada-idCAKBN2841E8"
$ws.Range("A16").Value = "Canada plans digital tax in 2022 on global tech giants - Reuters Canada"
$ws.Range("B16").Value = "https://ca.reuters.com/article/us-canada-budget-tax-idCAKBN28A2ZM"
$ws.Range("A17").Value = "Bank of Canada reiterates it could cut rates further if COVID worsens - Reuters Canada"
$ws.Range("B17").Value = "https://ca.reuters.com/article/us-canada-cenbank-idCAKBN28K2WW"
$ws.Range("A18").Value = "Canada trade deficit shrinks slightly in October, still higher than expected - Reuters Canada"
$ws.Range("B18").Value = "https://ca.reuters.com/article/us-canada-economy-trade-idCAKBN28E2ZA"
$ws.Range("A19").Value = "CANADA FX DEBT-C`$ climbs as market shrugs off Bank of Canada jawboning - Reuters India"
$ws.Range("B19").Value = "https://uk.reuters.com/article/canada-forex-idINL1N2IV2FS"
$ws.Range("A20").Value = "Air Canada makes more cuts in Atlantic Canada - CBC News"
$ws.Range("B20").Value = "https://www.youtube.com/watch?v=m9AwKeuVfeU"
$ws.Range("A21").Value = "Canadian Fashion Mogul Peter Nygard Indicted on Sex-Trafficking Charges"
$ws.Range("B21").Value = "https://www.nytimes.com/2020/12/15/world/canada/peter-nygard-sex-trafficking-charges.html"
